# "Add files via upload" — re-uploaded competitive-analysis sheet with a few
# manual tweaks made in Excel before saving:
#   1. C4: "詹惟中面相大師-面相算命必備(簡稱A)" -> drop the trailing "(簡稱A)"
#   2. D4: "面相大師—AI掃臉看面相與手相App(簡稱B)" -> drop the trailing "(簡稱B)"
#   3. D9 (下載量 / download count for App B): was blank -> now "_"
#      (matches the "_" placeholder already used at D7 for "not applicable")
#   4. Leave the selection on D13 afterwards (where the user's cursor ended up)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) C4: remove the "(簡稱A)" suffix (5 chars), same as a user selecting
#          those characters and pressing Delete -- leaves the remaining runs
#          ("詹惟中面相大師" / "-" / "面相算命必備") and their formatting intact. ---
$c4 = $ws.Range("C4")
$c4Text = $c4.Value2
$c4.Characters($c4Text.Length - 4, 5).Text = ""

# re-assert the per-run formatting that belonged to the text that remains,
# in case the character-delete re-flowed the run boundaries
$ws.Range("C4").Characters(8, 1).Font.Name = "Times New Roman"
$ws.Range("C4").Characters(8, 1).Font.Size = 14
$ws.Range("C4").Characters(9, 6).Font.Name = "標楷體"
$ws.Range("C4").Characters(9, 6).Font.Size = 14

# --- 2) D4: remove the "(簡稱B)" suffix (5 chars) the same way ---
$d4 = $ws.Range("D4")
$d4Text = $d4.Value2
$d4.Characters($d4Text.Length - 4, 5).Text = ""

$ws.Range("D4").Characters(6, 2).Font.Name = "Segoe UI Historic"
$ws.Range("D4").Characters(6, 2).Font.Size = 14
$ws.Range("D4").Characters(8, 8).Font.Name = "標楷體"
$ws.Range("D4").Characters(8, 8).Font.Size = 14
$ws.Range("D4").Characters(16, 3).Font.Name = "Segoe UI Historic"
$ws.Range("D4").Characters(16, 3).Font.Size = 14

# --- 3) D9: fill in the missing download-count cell with "_" ---
$ws.Range("D9").Value = "_"

# --- 4) match the final cursor position recorded in the saved file ---
$ws.Range("D13").Select()
